$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (order matters for shared-string table ordering) ---
$ws.Range("A1").Value = "Cod_Articolo"
$ws.Range("B1").Value = "Barcode"
$ws.Range("C1").Value = "Descri_Articolo"
$ws.Range("D1").Value = "Pz_x_Conf"
$ws.Range("E1").Value = "Prezzo"
$ws.Range("F1").Value = "Famiglia"

# --- Article data, row by row, left to right ---
$data = @(
  @(1,  111111, "DHSHSHSDHSD",   10, 9.99,  10),
  @(2,  222222, "SYSEYTYETYRTY", 10, 56.2,  15),
  @(3,  333333, "WYRTHFGHFG",    10, 14.5,  16),
  @(4,  444444, "HSTHDDJHRUY",   10, 45,    17),
  @(5,  555555, "SDRGSDFGSDF",   10, 47.1,  20),
  @(6,  666666, "GSDFGSDFG",     10, 56.8,  25),
  @(7,  777777, "SDGSDFGSDFG",   10, 85,    28),
  @(8,  888888, "SDFGSDFGSDFG",  10, 2.36,  29),
  @(9,  999999, "SDFGSDFGSDFG",  10, 1.87,  32),
  @(10, 101010, "SDFGSDGS",      10, 55.23, 35)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 1).HorizontalAlignment = -4131
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $r++
}

# --- Column widths (auto-fit like Excel would do after data entry) ---
$ws.Range("A1:D1").EntireColumn.AutoFit()

# --- Selection like in the saved file ---
[void]$ws.Range("M10").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
